# Appends six new daily rows (2025-12-16 through 2025-12-18, two stations
# each) to Sheet1, mirroring the layout/formulas of the existing rows, then
# updates the sheet's active selection to the new bottom cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Clone formatting (styles) of the last existing data row (213) down
#        into the six new rows, one at a time so each destination row gets
#        the same per-column style indices as the source row.
for ($r = 214; $r -le 219; $r++) {
    $ws.Range("A213:H213").Copy()
    $ws.Range("A$r`:H$r").PasteSpecial(-4104)
}
$excel.CutCopyMode = $false

# --- 2. Row 214: 2025-12-16, 四方坪站
$ws.Range("A214").Value = 46007
$ws.Range("B214").Value = "四方坪站"
$ws.Range("C214").Formula = "=15455/126"
$ws.Range("D214").Formula = "=C214/1440"
$ws.Range("E214").Formula = "=7975.08/126"
$ws.Range("F214").Formula = "=2617.95/126"
$ws.Range("G214").Formula = "=7975.08/(15455/60)"
$ws.Range("H214").Formula = "=345/126"

# --- 3. Row 215: 2025-12-16, 高岭站
$ws.Range("A215").Value = 46007
$ws.Range("B215").Value = "高岭站"
$ws.Range("C215").Formula = "=6397/36"
$ws.Range("D215").Formula = "=C215/1440"
$ws.Range("E215").Formula = "=4366.49/36"
$ws.Range("F215").Formula = "=1160.77/36"
$ws.Range("G215").Formula = "=4366.49/(6397/60)"
$ws.Range("H215").Formula = "=164/36"

# --- 4. Row 216: 2025-12-17, 四方坪站
$ws.Range("A216").Value = 46008
$ws.Range("B216").Value = "四方坪站"
$ws.Range("C216").Formula = "=16475/126"
$ws.Range("D216").Formula = "=C216/1440"
$ws.Range("E216").Formula = "=7802.58/126"
$ws.Range("F216").Formula = "=2592.82/126"
$ws.Range("G216").Formula = "=7802/(15475/60)"
$ws.Range("H216").Formula = "=349/126"

# --- 5. Row 217: 2025-12-17, 高岭站
$ws.Range("A217").Value = 46008
$ws.Range("B217").Value = "高岭站"
$ws.Range("C217").Formula = "=7017/36"
$ws.Range("D217").Formula = "=C217/1440"
$ws.Range("E217").Formula = "=4438.31/36"
$ws.Range("F217").Formula = "=1209.84/36"
$ws.Range("G217").Formula = "=4436.31/(7017/60)"
$ws.Range("H217").Formula = "=169/36"

# --- 6. Row 218: 2025-12-18, 四方坪站
$ws.Range("A218").Value = 46009
$ws.Range("B218").Value = "四方坪站"
$ws.Range("C218").Formula = "=16700/126"
$ws.Range("D218").Formula = "=C218/1440"
$ws.Range("E218").Formula = "=8362.59/126"
$ws.Range("F218").Formula = "=2767.43/126"
$ws.Range("G218").Formula = "=8362.59/(16700/60)"
$ws.Range("H218").Formula = "=374/126"

# --- 7. Row 219: 2025-12-18, 高岭站
$ws.Range("A219").Value = 46009
$ws.Range("B219").Value = "高岭站"
$ws.Range("C219").Formula = "=6960/36"
$ws.Range("D219").Formula = "=C219/1440"
$ws.Range("E219").Formula = "=4338.56/36"
$ws.Range("F219").Formula = "=1117.58/36"
$ws.Range("G219").Formula = "=4338.56/(6960/60)"
$ws.Range("H219").Formula = "=174/36"

# --- 8. Move the visible selection to the new bottom of the sheet (the
#        author was last positioned at I220 after entering the data).
[void]$ws.Range("I220").Select()
